$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09102366666666667
$ws.Range("H2").Value = 0.273071
$ws.Range("I2").Value = 0.001049496534254756
$ws.Range("J2").Value = 0.001049496534254756
$ws.Range("M2").Value = 0.428743
$ws.Range("N2").Value = 1.286229
$ws.Range("O2").Value = 0.00412050394863168
$ws.Range("P2").Value = 0.00412050394863168
$ws.Range("Q2").Value = 0.03902575991766667
$ws.Range("R2").Value = 0.351231839259
$ws.Range("S2").Value = 0.000004324454613471987
$ws.Range("T2").Value = 0.000004324454613471987

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09102366666666667
$ws.Range("H3").Value = 0.273071
$ws.Range("I3").Value = 0.001049496534254756
$ws.Range("J3").Value = 0.001049496534254756
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("O3").Value = 0.7710272268990069
$ws.Range("P3").Value = 0.7710272268990069
$ws.Range("Q3").Value = 7.302486254608999
$ws.Range("R3").Value = 65.722376291481
$ws.Range("S3").Value = 0.0008091904024465633
$ws.Range("T3").Value = 0.0008091904024465633

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09102366666666667
$ws.Range("H4").Value = 0.273071
$ws.Range("I4").Value = 0.001049496534254756
$ws.Range("J4").Value = 0.001049496534254756
$ws.Range("M4").Value = 23.39612766666667
$ws.Range("N4").Value = 70.188383
$ws.Range("O4").Value = 0.2248522691523614
$ws.Range("P4").Value = 0.2248522691523614
$ws.Range("Q4").Value = 2.129601326021445
$ws.Range("R4").Value = 19.166411934193
$ws.Range("S4").Value = 0.0002359816771947209
$ws.Range("T4").Value = 0.0002359816771947209

# Row 5
$ws.Range("G5").Value = 83.68760400000001
$ws.Range("I5").Value = 0.9649122428754918
$ws.Range("J5").Value = 0.9649122428754917
$ws.Range("M5").Value = 0.428743
$ws.Range("N5").Value = 1.286229
$ws.Range("O5").Value = 0.00412050394863168
$ws.Range("P5").Value = 0.00412050394863168
$ws.Range("Q5").Value = 35.88047440177201
$ws.Range("R5").Value = 322.924269615948
$ws.Range("S5").Value = 0.003975924706851515
$ws.Range("T5").Value = 0.003975924706851514

# Row 6
$ws.Range("G6").Value = 83.68760400000001
$ws.Range("I6").Value = 0.9649122428754918
$ws.Range("J6").Value = 0.9649122428754917
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("O6").Value = 0.7710272268990069
$ws.Range("P6").Value = 0.7710272268990069
$ws.Range("Q6").Value = 6713.941552466147
$ws.Range("R6").Value = 60425.47397219532
$ws.Range("S6").Value = 0.7439736108251915
$ws.Range("T6").Value = 0.7439736108251914

# Row 7
$ws.Range("G7").Value = 83.68760400000001
$ws.Range("I7").Value = 0.9649122428754918
$ws.Range("J7").Value = 0.9649122428754917
$ws.Range("M7").Value = 23.39612766666667
$ws.Range("N7").Value = 70.188383
$ws.Range("O7").Value = 0.2248522691523614
$ws.Range("P7").Value = 0.2248522691523614
$ws.Range("Q7").Value = 1957.965867301444
$ws.Range("R7").Value = 17621.692805713
$ws.Range("S7").Value = 0.2169627073434488
$ws.Range("T7").Value = 0.2169627073434488

# Row 8
$ws.Range("G8").Value = 2.952165333333333
$ws.Range("H8").Value = 8.856496
$ws.Range("I8").Value = 0.03403826059025349
$ws.Range("J8").Value = 0.03403826059025349
$ws.Range("M8").Value = 0.428743
$ws.Range("N8").Value = 1.286229
$ws.Range("O8").Value = 0.00412050394863168
$ws.Range("P8").Value = 0.00412050394863168
$ws.Range("Q8").Value = 1.265720221509333
$ws.Range("R8").Value = 11.391481993584
$ws.Range("S8").Value = 0.0001402547871666936
$ws.Range("T8").Value = 0.0001402547871666936

# Row 9
$ws.Range("G9").Value = 2.952165333333333
$ws.Range("H9").Value = 8.856496
$ws.Range("I9").Value = 0.03403826059025349
$ws.Range("J9").Value = 0.03403826059025349
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("O9").Value = 0.7710272268990069
$ws.Range("P9").Value = 0.7710272268990069
$ws.Range("Q9").Value = 236.841115695184
$ws.Range("R9").Value = 2131.570041256656
$ws.Range("S9").Value = 0.0262444256713689
$ws.Range("T9").Value = 0.0262444256713689

# Row 10
$ws.Range("G10").Value = 2.952165333333333
$ws.Range("H10").Value = 8.856496
$ws.Range("I10").Value = 0.03403826059025349
$ws.Range("J10").Value = 0.03403826059025349
$ws.Range("M10").Value = 23.39612766666667
$ws.Range("N10").Value = 70.188383
$ws.Range("O10").Value = 0.2248522691523614
$ws.Range("P10").Value = 0.2248522691523614
$ws.Range("Q10").Value = 69.06923703177422
$ws.Range("R10").Value = 621.623133285968
$ws.Range("S10").Value = 0.007653580131717894
$ws.Range("T10").Value = 0.007653580131717894
